$d = $word.ActiveDocument

# Namespace wrapper reused for every InsertXML payload below. Each call
# rebuilds a whole <w:p> (keeping the original paragraph's attributes and
# <w:pPr>) but with a fresh set of <w:r> runs, which is the only reliable
# way to get distinct (non-merged) sibling runs out of this host -- plain
# Range.Text / Range.InsertAfter silently collapse same-formatted adjacent
# runs back into one on save.
#
# NOTE: InsertXML-ing a paragraph that is currently the *last* paragraph in
# the body leaves behind a stray trailing empty paragraph, so every such
# replacement below happens while the paragraph still has a following
# sibling; the paragraph-7 deletion (which shrinks the body) is done last.
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-ParagraphXml {
    param($Paragraph, $InnerXml)
    $r = $Paragraph.Range
    $xml = $pkgOpen + $InnerXml + $pkgClose
    $r.InsertXML($xml)
}

# --- Paragraph 2: "DAY 1" -> "DAY " + "2" ---
$p2 = $d.Paragraphs(2)
$inner2 = '<w:p w14:paraId="7D48947F" w14:textId="6AC3C0A1" w:rsidR="00DA763F" w:rsidRDefault="006512A5" w:rsidP="006512A5"><w:pPr><w:pStyle w:val="Title"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">DAY </w:t></w:r><w:r><w:t>2</w:t></w:r></w:p>'
Replace-ParagraphXml $p2 $inner2

# --- Paragraph 4: "Created reusable button component" -> "Completed the Services section " ---
$p4 = $d.Paragraphs(4)
$inner4 = '<w:p w14:paraId="66913D24" w14:textId="22E6D888" w:rsidR="006512A5" w:rsidRDefault="006512A5" w:rsidP="006512A5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Completed the Services section </w:t></w:r></w:p>'
Replace-ParagraphXml $p4 $inner4

# --- Paragraph 5: "Created Nav Bar" -> "Completed" + " Reasons" + " Section" ---
$p5 = $d.Paragraphs(5)
$inner5 = '<w:p w14:paraId="5C69A5B2" w14:textId="0DB74FFF" w:rsidR="006512A5" w:rsidRDefault="006512A5" w:rsidP="006512A5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Completed</w:t></w:r><w:r><w:t xml:space="preserve"> Reasons</w:t></w:r><w:r><w:t xml:space="preserve"> Section</w:t></w:r></w:p>'
Replace-ParagraphXml $p5 $inner5

# --- Paragraph 6: "Completed Hero Section" -> "Completed " + "CTA" + " section" ---
$p6 = $d.Paragraphs(6)
$inner6 = '<w:p w14:paraId="19A2AFA5" w14:textId="31940693" w:rsidR="006512A5" w:rsidRDefault="006512A5" w:rsidP="006512A5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Completed </w:t></w:r><w:r><w:t>CTA</w:t></w:r><w:r><w:t xml:space="preserve"> section</w:t></w:r></w:p>'
Replace-ParagraphXml $p6 $inner6

# --- Paragraph 7: "Completed Info section" -> delete the whole paragraph ---
# Delete from just before paragraph 6's own paragraph mark through the end
# of paragraph 7; this removes paragraph 7's run + pPr + its own mark and
# merges what's left into paragraph 6 (which keeps paragraph 6's mark/pPr).
$p6 = $d.Paragraphs(6)
$p7 = $d.Paragraphs(7)
$mergeRange = $d.Range($p6.Range.End - 1, $p7.Range.End)
$mergeRange.Delete() | Out-Null

Write-Host "Done. Final paragraphs:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Host $i ":" $p.Range.Text
}
